$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45171 -> 45172) for every data row (rows 2 through 130).
$ws.Range("C2:C130").Value = 45172
